# Add a new header/question row at the top of the sheet (rows shift down
# conceptually is NOT what happens here - instead row 1/2 contents are
# replaced with new "product"/"text" style header values, and new cells
# are added in columns B and C for rows 1 and 2).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "product1"
$ws.Range("B1").Value = "prodcut2"
$ws.Range("C1").Value = "product3"

$ws.Range("A2").Value = "text"
$ws.Range("B2").Value = "radio"
$ws.Range("C2").Value = "checkbox"

[void]$ws.Range("C2").Select()
